$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in column C for rows 2 through 32
# from 45605 to 45606 (one day later).
for ($r = 2; $r -le 32; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45605) {
        $cell.Value = 45606
    }
}

# Swap the A (Beteckning) and G (Area (ha)) values between row 30 and row 31.
$a30 = $ws.Cells.Item(30, 1).Value2
$a31 = $ws.Cells.Item(31, 1).Value2
$ws.Cells.Item(30, 1).Value = $a31
$ws.Cells.Item(31, 1).Value = $a30

$g30 = $ws.Cells.Item(30, 7).Value2
$g31 = $ws.Cells.Item(31, 7).Value2
$ws.Cells.Item(30, 7).Value = $g31
$ws.Cells.Item(31, 7).Value = $g30
